$d = $word.ActiveDocument

$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx += 1
    $t = $p.Range.Text
    if ($t -like "*train set*") {
        $target = $p
        $targetIdx = $idx
    }
}

$target.Range.InsertParagraphAfter() | Out-Null

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx += 1
    if ($idx -eq $targetIdx + 1) {
        $newPara = $p
    }
}

$newPara.Range.Text = "We will calculate the metrics to predict the accuracy of the predicted data and the predicted data is visualized using the matplotlib."

# re-fetch fresh paragraph reference
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx += 1
    if ($idx -eq $targetIdx + 1) {
        $newParaFresh = $p
    }
}
$pStart = $newParaFresh.Range.Start
$pEnd = $newParaFresh.Range.End
Write-Output "pStart=$pStart pEnd=$pEnd"

$bmStart = $pEnd - 1
Write-Output "bmStart=$bmStart"
$bmRange = $d.Range($bmStart, $bmStart)
Write-Output "bmRange.Start=$($bmRange.Start) End=$($bmRange.End)"

$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
Write-Output "bookmark added"
